# Update "想去人数" (column F) values across the relevant worksheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1601
$ws1.Range("F4").Value  = 8575
$ws1.Range("F6").Value  = 71
$ws1.Range("F7").Value  = 26
$ws1.Range("F8").Value  = 83
$ws1.Range("F9").Value  = 1350
$ws1.Range("F10").Value = 114
$ws1.Range("F11").Value = 25
$ws1.Range("F12").Value = 29
$ws1.Range("F13").Value = 9255
$ws1.Range("F15").Value = 92
$ws1.Range("F16").Value = 215
$ws1.Range("F17").Value = 172
$ws1.Range("F18").Value = 348
$ws1.Range("F19").Value = 6157
$ws1.Range("F20").Value = 1052
$ws1.Range("F21").Value = 68
$ws1.Range("F22").Value = 39
$ws1.Range("F23").Value = 110

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 33

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1601
$ws4.Range("F4").Value  = 8575
$ws4.Range("F6").Value  = 71
$ws4.Range("F7").Value  = 26
$ws4.Range("F8").Value  = 83
$ws4.Range("F9").Value  = 1350
$ws4.Range("F10").Value = 114
$ws4.Range("F11").Value = 25
$ws4.Range("F12").Value = 29
$ws4.Range("F13").Value = 33
$ws4.Range("F15").Value = 9255
$ws4.Range("F17").Value = 92
$ws4.Range("F18").Value = 215
$ws4.Range("F19").Value = 172
$ws4.Range("F20").Value = 348
$ws4.Range("F21").Value = 6157
$ws4.Range("F22").Value = 1052
$ws4.Range("F23").Value = 68
$ws4.Range("F24").Value = 39
$ws4.Range("F25").Value = 110

$wb.Save()
